# Apply updated cryptocurrency price/volume data (and 3 re-ordered rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.055.02'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.208.15'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.54'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.39'
$ws.Range('E6').Value = '  -6.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.199.53'
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('E10').Value = '  -4.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.23'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.475'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000231'
$ws.Range('E13').Value = '  -3.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.58'
$ws.Range('E14').Value = '  -7.20%  '
$ws.Range('D15').Value = '3.715.76'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').Value = '66.990.01'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '3.209.60'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.80'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '499.64'
$ws.Range('E20').Value = '  -3.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.16'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.710'
$ws.Range('E22').Value = '  -4.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.31'
$ws.Range('E23').Value = '  -4.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.25'
$ws.Range('E24').Value = '  -3.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.80'
$ws.Range('E25').Value = '  -2.77%  '
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '27.60'
$ws.Range('E28').Value = '  -3.02%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.02'
$ws.Range('E29').Value = '  -5.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.55'
$ws.Range('E30').Value = '  -4.88%  '
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.51'
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '516.97'
$ws.Range('E33').Value = '  -2.15%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '54.31'
$ws.Range('E35').Value = '  +1.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.02'
$ws.Range('E36').Value = '  -5.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.22'
$ws.Range('E37').Value = '  -7.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0413'
$ws.Range('E38').Value = '  -3.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0807'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.54'
$ws.Range('E40').Value = '  -6.06%  '
$ws.Range('E41').Value = '  -3.52%  '
$ws.Range('D42').Value = '2.847.81'
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.47'
$ws.Range('E43').Value = '  -11.06%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.248'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '121.13'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.54'
$ws.Range('E47').Value = '  -4.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.00'
$ws.Range('E48').Value = '  -6.26%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.108'
$ws.Range('E49').Value = '  -2.63%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0519'
$ws.Range('E50').Value = '  -10.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.09'
$ws.Range('E51').Value = '  -12.21%  '
